$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells whose new values would otherwise
# be reinterpreted as numbers by Excel (losing trailing zeros / becoming
# scientific notation), so the exact textual content from the source is kept.
$ws.Range("D4:D11").NumberFormat = "@"
$ws.Range("D13:D20").NumberFormat = "@"
$ws.Range("D22:D24").NumberFormat = "@"
$ws.Range("D26:D51").NumberFormat = "@"

# --- Update Price / Volume(1h) columns for rows 2-33 ---
$ws.Cells.Item(2, 4).Value = "29.541.08"
$ws.Cells.Item(2, 5).Value = "  +0.30%  "
$ws.Cells.Item(3, 4).Value = "1.905.28"
$ws.Cells.Item(3, 5).Value = "  -0.65%  "
$ws.Cells.Item(4, 4).Value = "1.013"
$ws.Cells.Item(4, 5).Value = "  +0.44%  "
$ws.Cells.Item(5, 4).Value = "338.05"
$ws.Cells.Item(5, 5).Value = "  +3.95%  "
$ws.Cells.Item(6, 4).Value = "1.013"
$ws.Cells.Item(6, 5).Value = "  +0.53%  "
$ws.Cells.Item(7, 4).Value = "0.4769"
$ws.Cells.Item(7, 5).Value = "  -1.26%  "
$ws.Cells.Item(8, 4).Value = "0.4010"
$ws.Cells.Item(8, 5).Value = "  -1.68%  "
$ws.Cells.Item(9, 4).Value = "0.08050"
$ws.Cells.Item(9, 5).Value = "  -2.23%  "
$ws.Cells.Item(10, 4).Value = "0.9896"
$ws.Cells.Item(10, 5).Value = "  -2.62%  "
$ws.Cells.Item(11, 4).Value = "23.29"
$ws.Cells.Item(11, 5).Value = "  -0.55%  "
$ws.Cells.Item(12, 4).Value = "1.926.72"
$ws.Cells.Item(12, 5).Value = "  +1.46%  "
$ws.Cells.Item(13, 4).Value = "5.930"
$ws.Cells.Item(13, 5).Value = "  -2.60%  "
$ws.Cells.Item(14, 4).Value = "7.120"
$ws.Cells.Item(14, 5).Value = "  -1.72%  "
$ws.Cells.Item(15, 4).Value = "89.40"
$ws.Cells.Item(15, 5).Value = "  -2.21%  "
$ws.Cells.Item(16, 4).Value = "0.06840"
$ws.Cells.Item(16, 5).Value = "  +0.19%  "
$ws.Cells.Item(17, 4).Value = "1.013"
$ws.Cells.Item(17, 5).Value = "  +0.46%  "
$ws.Cells.Item(18, 4).Value = "0.00001021"
$ws.Cells.Item(18, 5).Value = "  -1.83%  "
$ws.Cells.Item(19, 4).Value = "17.37"
$ws.Cells.Item(19, 5).Value = "  -2.01%  "
$ws.Cells.Item(20, 4).Value = "1.011"
$ws.Cells.Item(20, 5).Value = "  +0.41%  "
$ws.Cells.Item(21, 4).Value = "29.568.98"
$ws.Cells.Item(21, 5).Value = "  +0.34%  "
$ws.Cells.Item(22, 4).Value = "5.519"
$ws.Cells.Item(22, 5).Value = "  -2.44%  "
$ws.Cells.Item(23, 4).Value = "11.63"
$ws.Cells.Item(23, 5).Value = "  -1.63%  "
$ws.Cells.Item(24, 4).Value = "2.152"
$ws.Cells.Item(24, 5).Value = "  -1.23%  "
$ws.Cells.Item(25, 4).Value = "2.181.88"
$ws.Cells.Item(25, 5).Value = "  +1.76%  "
$ws.Cells.Item(26, 4).Value = "156.82"
$ws.Cells.Item(26, 5).Value = "  +0.43%  "
$ws.Cells.Item(27, 4).Value = "6.461"
$ws.Cells.Item(27, 5).Value = "  -2.48%  "
$ws.Cells.Item(28, 4).Value = "19.72"
$ws.Cells.Item(28, 5).Value = "  -1.88%  "
$ws.Cells.Item(29, 4).Value = "2.053"
$ws.Cells.Item(29, 5).Value = "  -3.03%  "
$ws.Cells.Item(30, 4).Value = "119.42"
$ws.Cells.Item(30, 5).Value = "  -1.00%  "
$ws.Cells.Item(31, 4).Value = "0.9943"
$ws.Cells.Item(31, 5).Value = "  -2.87%  "
$ws.Cells.Item(32, 4).Value = "0.09526"
$ws.Cells.Item(32, 5).Value = "  -0.76%  "
$ws.Cells.Item(33, 4).Value = "5.488"
$ws.Cells.Item(33, 5).Value = "  -3.48%  "

# --- Rows 34 and 35 swapped places (coin identity + data) ---
$ws.Cells.Item(34, 2).Value = "HuobiToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(34, 4).Value = "3.544"
$ws.Cells.Item(34, 5).Value = "  -0.29%  "
$ws.Cells.Item(35, 2).Value = "ARBITRUM"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(35, 4).Value = "1.385"
$ws.Cells.Item(35, 5).Value = "  +0.91%  "

# --- Update Price / Volume(1h) columns for rows 36-51 ---
$ws.Cells.Item(36, 4).Value = "0.06479"
$ws.Cells.Item(36, 5).Value = "  +5.94%  "
$ws.Cells.Item(37, 4).Value = "0.02240"
$ws.Cells.Item(37, 5).Value = "  -2.21%  "
$ws.Cells.Item(38, 4).Value = "1.195"
$ws.Cells.Item(38, 5).Value = "  +0.94%  "
$ws.Cells.Item(39, 4).Value = "0.5828"
$ws.Cells.Item(39, 5).Value = "  -2.92%  "
$ws.Cells.Item(40, 4).Value = "10.57"
$ws.Cells.Item(40, 5).Value = "  -2.33%  "
$ws.Cells.Item(41, 4).Value = "7.788"
$ws.Cells.Item(41, 5).Value = "  -3.99%  "
$ws.Cells.Item(42, 4).Value = "0.1820"
$ws.Cells.Item(42, 5).Value = "  -1.68%  "
$ws.Cells.Item(43, 4).Value = "2.455"
$ws.Cells.Item(43, 5).Value = "  +1.65%  "
$ws.Cells.Item(44, 4).Value = "1.241"
$ws.Cells.Item(44, 5).Value = "  -2.97%  "
$ws.Cells.Item(45, 4).Value = "0.07428"
$ws.Cells.Item(45, 5).Value = "  -2.44%  "
$ws.Cells.Item(46, 4).Value = "12.09"
$ws.Cells.Item(46, 5).Value = "  -3.25%  "
$ws.Cells.Item(47, 4).Value = "0.5468"
$ws.Cells.Item(47, 5).Value = "  -2.29%  "
$ws.Cells.Item(48, 4).Value = "1.941"
$ws.Cells.Item(48, 5).Value = "  -0.88%  "
$ws.Cells.Item(49, 4).Value = "116.02"
$ws.Cells.Item(49, 5).Value = "  -2.00%  "
$ws.Cells.Item(50, 4).Value = "2.384"
$ws.Cells.Item(50, 5).Value = "  -1.89%  "
$ws.Cells.Item(51, 4).Value = "71.27"
$ws.Cells.Item(51, 5).Value = "  -1.70%  "
